# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 4250
$wsExhibition.Range("F3").Value = 2418
$wsExhibition.Range("F10").Value = 129
$wsExhibition.Range("F12").Value = 1584
$wsExhibition.Range("F14").Value = 3264

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 4250
$wsAll.Range("F3").Value = 2418
$wsAll.Range("F12").Value = 129
$wsAll.Range("F16").Value = 1584
$wsAll.Range("F18").Value = 3264
